# Update Column D ("Final Value") with specific fitness calculation values
# based on Nutrient signal (issue #56).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 6141.95
    3  = 992.47
    4  = 179.42
    5  = 182.88
    6  = 81.88
    7  = 257.6
    8  = 1379.76
    9  = 3.2
    10 = 48.93
    11 = 51.65
    12 = 5467.4
    13 = 5674.82
    14 = 27.5
    15 = 751.48
    16 = 1.99
    17 = 1.65
    18 = 2.34
    19 = 29.44
    20 = 4.59
    21 = 75.26000000000001
    22 = 602.02
    23 = 3.15
    24 = 1289.85
    25 = 2413.84
    26 = 31.69
}

foreach ($row in $newValues.Keys) {
    $ws.Range("D$row").Value = $newValues[$row]
}
